$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Qminus1)
$ws.Range("B2").Value = 0.008240683130927028
$ws.Range("C2").Value = 0.3653326498374267
$ws.Range("D2").Value = 0.2391193703417652
$ws.Range("E2").Value = 0.4889983336799474
$ws.Range("F2").Value = 0.507385489981291

# Row 3 (Q0)
$ws.Range("B3").Value = 0.05900012108174069
$ws.Range("C3").Value = 0.3849571187190655
$ws.Range("D3").Value = 0.2374917410917716
$ws.Range("E3").Value = 0.4873312437057279
$ws.Range("F3").Value = 0.5007252527271727
$ws.Range("G3").Value = 15

# Row 4 (Q1)
$ws.Range("B4").Value = 0.3329153595043876
$ws.Range("C4").Value = 0.56922599019722
$ws.Range("D4").Value = 0.4729138842649011
$ws.Range("E4").Value = 0.6876873448485882
$ws.Range("F4").Value = 0.6244466761365319
$ws.Range("G4").Value = 14
